$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "May 2019" row (row 10) below the existing "April 2019" row.
$ws.Range("A10").Value = "May 2019"
$ws.Range("B10").Value = "https://myemail.constantcontact.com/News-From-The-Forest---May-2019.html?soid=1102494320279&aid=11jtFPoUhxc"

# Turn B10 into a real hyperlink, same as the other month rows above it.
$ws.Hyperlinks.Add($ws.Range("B10"), "https://myemail.constantcontact.com/News-From-The-Forest---May-2019.html?soid=1102494320279&aid=11jtFPoUhxc")

# New row's link cell gets vertically-centered text (new cell style).
$ws.Range("B10").VerticalAlignment = -4108

# Move the active selection down, matching the refreshed view state.
$ws.Range("B25").Select() | Out-Null
